$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.022520065307617
$ws.Range("B1").Value = 3.454584121704102
$ws.Range("C1").Value = 4.533708095550537
$ws.Range("D1").Value = 2.05327320098877
$ws.Range("E1").Value = 1.606835126876831
